$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.372.48"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "3.138.37"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.93"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.04"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "3.135.09"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.48"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "3.652.60"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("E16").Value = "  +2.57%  "

$ws.Range("D17").Value = "64.331.34"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "3.154.70"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.97"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.89"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.80"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.55"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.49"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -3.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.49"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +7.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.116"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +2.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -5.23%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.99"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.52%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.64"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.84%  "

$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.56"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("D38").Value = "0.0₃0748"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "447.15"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.31"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("D44").Value = "2.893.87"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +5.60%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.36"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.114"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.35"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.88%  "
